# AutoCommit_6 ноября 2023 г. 16:01:52_SibNout2023
# Adds a spacer row, reorders/extends the Лаб/Инд header columns, removes
# the now-unused G column data, and appends a new student row (Шаповаленко).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New student "Шаповаленко" appended as row 34 (adds the shared string
#    first so it lands at the same sharedStrings index as the target).
# ---------------------------------------------------------------------
$ws.Range("B34").Value = "Шаповаленко"
$ws.Range("B34").Font.Bold = $true
$ws.Range("B34").HorizontalAlignment = -4108
$ws.Range("B34").VerticalAlignment = -4108
$ws.Range("B34").WrapText = $true

$ws.Range("F34").Value = 5
$ws.Range("F34").HorizontalAlignment = -4108
$ws.Range("F34").VerticalAlignment = -4108
$ws.Range("F34").WrapText = $true

# ---------------------------------------------------------------------
# 2) Header row 3: Лаб1/Инд2/Инд3 shuffled, 3 new lab columns appended.
#    Инд1 (C3) is untouched.
# ---------------------------------------------------------------------
$ws.Range("D3").Value = "Инд2"
$ws.Range("E3").Value = "Инд3"
$ws.Range("F3").Value = "Лаб1"
$ws.Range("G3").Value = "Лаб2"
$ws.Range("H3").Value = "Лаб3-4"
$ws.Range("I3").Value = "Лаб5"

# ---------------------------------------------------------------------
# 3) Data grid: Инд2/Инд3/Лаб1 scores reshuffled per student, old Лаб1
#    (column G, now unused) fully cleared for every row.
# ---------------------------------------------------------------------
$ws.Range("D5").Value = 5
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 5

$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 5

$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 5

$ws.Range("D8").Value = 5
$ws.Range("E8").ClearContents()
$ws.Range("F8").Value = 5

$ws.Range("D9").Value = 5
$ws.Range("E9").ClearContents()
$ws.Range("F9").Value = 5

$ws.Range("D10").Value = 5
$ws.Range("E10").ClearContents()
$ws.Range("F10").Value = 5

$ws.Range("D11").ClearContents()
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 5

$ws.Range("D12").Value = 5
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()

$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("F13").Value = 5

$ws.Range("D14").ClearContents()
$ws.Range("E14").Value = 5
$ws.Range("F14").ClearContents()

$ws.Range("D15").Value = 5
$ws.Range("E15").ClearContents()
$ws.Range("F15").ClearContents()

$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()
$ws.Range("F16").Value = 5

$ws.Range("D17").Value = 5
$ws.Range("E17").ClearContents()
$ws.Range("F17").Value = 5

$ws.Range("D18").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("F18").Value = 5

$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()

$ws.Range("D20").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("F20").ClearContents()

$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 5

$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("F22").Value = 5

$ws.Range("D23").Value = 5
$ws.Range("E23").ClearContents()
$ws.Range("F23").Value = 5

$ws.Range("D24").ClearContents()
$ws.Range("E24").ClearContents()
$ws.Range("F24").Value = 5

$ws.Range("D25").Value = 5
$ws.Range("E25").ClearContents()
$ws.Range("F25").ClearContents()

$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 5
$ws.Range("F26").Value = 5

$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = 5

$ws.Range("D28").Value = 5
$ws.Range("E28").ClearContents()
$ws.Range("F28").Value = 5

$ws.Range("D29").Value = 5
$ws.Range("E29").ClearContents()
$ws.Range("F29").Value = 5

$ws.Range("D30").Value = 5
$ws.Range("E30").ClearContents()
$ws.Range("F30").Value = 5

$ws.Range("D31").ClearContents()
$ws.Range("E31").ClearContents()
$ws.Range("F31").Value = 5

$ws.Range("D32").Value = 5
$ws.Range("E32").Value = 5
$ws.Range("F32").Value = 5

# Column G (old Лаб1 position) is no longer used anywhere in the grid.
$ws.Range("G3:G33").Clear()
$ws.Range("D4").Clear()

# ---------------------------------------------------------------------
# 4) New spacer row (row 2, previously a gap in the sheet) + new ht.
# ---------------------------------------------------------------------
$ws.Rows(2).RowHeight = 13

# ---------------------------------------------------------------------
# 5) Dimension / view bookkeeping: scroll to show the new columns and
#    select the first new lab header cell, matching the saved view.
# ---------------------------------------------------------------------
$ws.Range("G4").Select()
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 3
$aw.ScrollRow = 5
